{"js": "// Insert a new \"Otras decisiones\" list item right after the paragraph about\n// first publication being free, introducing the rule about auctions that\n// expire with no bids.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Los usuarios creados despu\u00e9s de la migraci\u00f3n poseen su primera publicaci\u00f3n gratis\";\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(anchorText) !== -1) {\n    anchorParagraph = p;\n    break;\n  }\n}\nif (!anchorParagraph) {\n  throw new Error(\"Could not locate the anchor paragraph about free first publication.\");\n}\n\n// Capture the list this bullet belongs to so the new paragraph joins the\n// exact same numbering definition (numId=1, ilvl=0) instead of minting one.\nconst list = anchorParagraph.list;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\nconst newParagraphText =\n  \"Aquellas subastas que venzan sin ninguna oferta realizada en la misma se dar\u00e1n por finalizada sin generar ninguna factura de compra ya que nadie adquiri\u00f3 el producto ofertado.\";\n\nconst newParagraph = anchorParagraph.insertParagraph(newParagraphText, Word.InsertLocation.after);\nnewParagraph.styleBuiltIn = Word.BuiltInStyleName.listParagraph;\nnewParagraph.attachToList(listId, 0);\n\nawait context.sync();\n", "ps1": "# Insert a new \"Otras decisiones\" bullet right after the paragraph about the\n# first publication being free, adding the rule about auctions that expire\n# without any bid.\n$d = $word.ActiveDocument\n\n$anchorText = \"Los usuarios creados despu\u00e9s de la migraci\u00f3n poseen su primera publicaci\u00f3n gratis\"\n$newText = \"Aquellas subastas que venzan sin ninguna oferta realizada en la misma se dar\u00e1n por finalizada sin generar ninguna factura de compra ya que nadie adquiri\u00f3 el producto ofertado.\"\n\n# Locate the paragraph that ends the existing bullet list item so the new\n# bullet is inserted immediately after it (and before the trailing blank\n# paragraph).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the anchor paragraph about free first publication.\"\n}\n\n$anchorParagraph = $d.Paragraphs.Item($targetIndex)\n\n# InsertParagraphAfter on a list-item Range clones that paragraph's pPr\n# (pStyle \"Prrafodelista\" + numPr ilvl=0/numId=1), so the new paragraph joins\n# the same bulleted list automatically.\n$anchorParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($targetIndex + 1)\n$newParagraph.Range.Text = $newText\n"}
